# ver 0.213 dodata opcija 'dodaj trosak'
# Adds a backup copy of the Users_input sheet (with all 15 original rows plus
# the 6 new user rows appended), trims Users_input itself down to just the 6
# new rows, and appends a new company record to Company_input.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Duplicate "Users_input" -> "Users_input (backup)", placed right after
#    the original sheet (this preserves all existing rows/hyperlinks/styles).
# ---------------------------------------------------------------------------
$usersInput = $wb.Worksheets.Item("Users_input")
$usersInput.Copy($null, $usersInput)
$backup = $wb.Worksheets.Item("Users_input (2)")
$backup.Name = "Users_input (backup)"

# Shared placeholder password hash re-used for every manually added account
# (same value already stored in Users_input!C13).
$pwdHash = $usersInput.Cells.Item(13, 3).Value()

# New user rows (ids 15-20) to append.
$newUsers = @(
    @{ Id = 15; Email = "steslo@gmail.com"; Name = "Slobodan"; Surname = "Stevanović"; Workplace = "Pomoćnik tesara" },
    @{ Id = 16; Email = "đordej@gmail.com"; Name = "Dejan "; Surname = "Đorđević"; Workplace = "Pomoćnik tesara" },
    @{ Id = 17; Email = "ljubomir.connect@gmail.com"; Name = "Ljubomir "; Surname = "Mitrović "; Workplace = "Pomoćni radnik" },
    @{ Id = 18; Email = "rajpen@gmail.com"; Name = "Pendyala "; Surname = "Rajashekhar"; Workplace = "pomoćni radnik" },
    @{ Id = 19; Email = "sampra@gmail.com"; Name = "Prathipati"; Surname = "Sampath"; Workplace = "Pomoćni radnik" },
    @{ Id = 20; Email = "grefer@gmail.com"; Name = "Fernandes "; Surname = "Greig Vincent"; Workplace = "Pomoćni radnik" }
)

# ---------------------------------------------------------------------------
# 2) Append the 6 new rows to the backup sheet (rows 16-21), keeping the
#    original 15 rows (1 header + 14 data rows... actually 1 header + 14 data
#    rows = 15 rows total) untouched.
# ---------------------------------------------------------------------------
$backupRow = $backup.UsedRange.Rows.Count + 1
foreach ($u in $newUsers) {
    $backup.Cells.Item($backupRow, 1).Value = $u.Id
    $backup.Cells.Item($backupRow, 2).Value = $u.Email
    $backup.Cells.Item($backupRow, 3).Value = $pwdHash
    $backup.Cells.Item($backupRow, 4).Value = $u.Name
    $backup.Cells.Item($backupRow, 5).Value = $u.Surname
    $backup.Cells.Item($backupRow, 6).Value = 1
    $backup.Cells.Item($backupRow, 7).Value = $u.Workplace
    $backup.Cells.Item($backupRow, 8).Value = "c_user"
    $backup.Cells.Item($backupRow, 9).Value = 4
    $backup.Cells.Item($backupRow, 10).Value = 0
    $backupRow = $backupRow + 1
}
$backup.Range("C16").Select()

# ---------------------------------------------------------------------------
# 3) Trim "Users_input" down to just the 6 new rows: drop the old hyperlinks
#    and the old data (rows 2-15), then write the new rows as rows 2-7.
# ---------------------------------------------------------------------------
$usersInput.Hyperlinks.Delete()
$usersInput.Range("A2:J15").Clear()

$row = 2
foreach ($u in $newUsers) {
    $usersInput.Cells.Item($row, 1).Value = $u.Id
    $usersInput.Cells.Item($row, 2).Value = $u.Email
    $usersInput.Cells.Item($row, 3).Value = $pwdHash
    $usersInput.Cells.Item($row, 4).Value = $u.Name
    $usersInput.Cells.Item($row, 5).Value = $u.Surname
    $usersInput.Cells.Item($row, 6).Value = 1
    $usersInput.Cells.Item($row, 7).Value = $u.Workplace
    $usersInput.Cells.Item($row, 8).Value = "c_user"
    $usersInput.Cells.Item($row, 9).Value = 4
    $usersInput.Cells.Item($row, 10).Value = 0
    $row = $row + 1
}

# ---------------------------------------------------------------------------
# 4) Append a new company record to "Company_input" (row 5).
# ---------------------------------------------------------------------------
$companyInput = $wb.Worksheets.Item("Company_input")
$companyInput.Cells.Item(5, 1).Value = 4
$companyInput.Cells.Item(5, 2).Value = "CONNECT GRADNJA 011 DOO"
$companyInput.Cells.Item(5, 3).Value = "ĆIRILA I METODIJA 2"
$companyInput.Cells.Item(5, 4).Value = 2
$companyInput.Cells.Item(5, 5).Value = 11000
$companyInput.Cells.Item(5, 6).Value = "BEOGRAD"
$companyInput.Cells.Item(5, 7).Value = "SRBIJA"
$companyInput.Cells.Item(5, 8).Value = 112762175
$companyInput.Cells.Item(5, 9).Value = 21734349
$companyInput.Cells.Item(5, 11).Value = "connectgradnja@gmail.com"
$companyInput.Cells.Item(5, 12).Value = 649194504
$companyInput.Cells.Item(5, 13).Value = "ef88c20087854a96.png"

$companyInput.Hyperlinks.Add($companyInput.Cells.Item(5, 14), "mailto:blagajnik3@gmail.com")
# match the hyperlink cell style used by the rows above (N2:N4)
$companyInput.Cells.Item(4, 14).Copy()
$companyInput.Cells.Item(5, 14).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$companyInput.Range("F15").Select()

# ---------------------------------------------------------------------------
# 5) Leave "Users_input" as the active/selected sheet, matching the source
#    workbook (tabSelected stays on this sheet).
# ---------------------------------------------------------------------------
$usersInput.Activate()
$usersInput.Range("C19").Select()
